$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update status of the remaining "Doing"/"To-Do" rows to "Done"
$ws.Range("C22").Value = "Done"
$ws.Range("C23").Value = "Done"
$ws.Range("C24").Value = "Done"
$ws.Range("C25").Value = "Done"
$ws.Range("C26").Value = "Done"

# Update selection / view state to match the author's final cursor position
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("C12").Select()
